$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new "Location" / "other location" / "Pincode" columns (E:G) ---
# Header row (row 1)
$ws.Range("E1").Value = "Location"
$ws.Range("F1").Value = "other location"
$ws.Range("G1").Value = "Pincode"

# Data row (row 2) - note: F2 (newyork) is written before E2 (Chennai) so the
# shared-string table fills up in the same order as the source workbook.
$ws.Range("F2").Value = "newyork"
$ws.Range("E2").Value = "Chennai"
$ws.Range("G2").Value = 1234567

# --- Resize columns E and F, add column F width ---
$ws.Columns.Item(5).ColumnWidth = 12.54296875
$ws.Columns.Item(6).ColumnWidth = 13.26953125

# --- Update the active selection shown in the sheet view ---
$ws.Range("C10").Select()
